$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new column D "canonical SMILES" (non-isomeric SMILES) next to the
# existing column C "canonical isomeric SMILES".
# ---------------------------------------------------------------------------

# Header
$ws.Range("D2").Value2 = "canonical SMILES"

# Row 3 - SM13_micro001 (has a distinct non-isomeric SMILES)
$ws.Range("D3").Value2 = "Cc1cccc(c1)[NH+]=c2c3cc(c(cc3[nH]cn2)OC)OC"

# Row 4 - SM13_micro003 (identical to column C, no stereo bonds present)
$ws.Range("D4").Value2 = "Cc1cccc(c1)[NH2+]c2c3cc(c(cc3[nH+]cn2)OC)OC"

# Row 5 - SM13_micro004
$ws.Range("D5").Value2 = "Cc1cccc(c1)[N-]c2c3cc(c(cc3ncn2)OC)OC"

# Row 6 - SM13_micro005
$ws.Range("D6").Value2 = "Cc1cccc(c1)Nc2c3cc(c(cc3ncn2)OC)OC"

# Row 7 - SM13_micro007 (has a distinct non-isomeric SMILES)
$ws.Range("D7").Value2 = "Cc1cccc(c1)N=c2c3cc(c(cc3nc[nH]2)OC)OC"

# Row 8 - SM13_micro008
$ws.Range("D8").Value2 = "Cc1cccc(c1)[NH2+]c2c3cc(c(cc3ncn2)OC)OC"

# Row 9 - SM13_micro009 (has a distinct non-isomeric SMILES)
$ws.Range("D9").Value2 = "Cc1cccc(c1)N=c2c3cc(c(cc3[nH]cn2)OC)OC"

# Row 10 - SM13_micro012
$ws.Range("D10").Value2 = "Cc1cccc(c1)Nc2c3cc(c(cc3nc[nH+]2)OC)OC"

# Row 11 - SM13_micro013
$ws.Range("D11").Value2 = "Cc1cccc(c1)Nc2c3cc(c(cc3[nH+]c[nH+]2)OC)OC"

# Row 12 - SM13_micro014
$ws.Range("D12").Value2 = "Cc1cccc(c1)[NH2+]c2c3cc(c(cc3nc[nH+]2)OC)OC"

# Row 13 - SM13_micro015
$ws.Range("D13").Value2 = "Cc1cccc(c1)[NH2+]c2c3cc(c(cc3[nH+]c[nH+]2)OC)OC"

# ---------------------------------------------------------------------------
# Copy the formatting (fill/font/border/alignment) from column C into the new
# column D, row by row, so the new cells reuse the same cell styles.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 13; $r++) {
    $src = $ws.Cells.Item($r, 3)
    $dst = $ws.Cells.Item($r, 4)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Widen the new column to fit its content (closest representable value to
# the target width of 37.7109375 given this engine's column-width quantization).
$ws.Columns.Item(4).ColumnWidth = 36.85
